# The "Förändrad" (Changed) date column C was bumped by one day
# (serial 46060 -> 46061, i.e. 2026-02-07 -> 2026-02-08) for every
# data row (rows 2 through 352).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C352").Value = 46061
